$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns A (dates) and C (amounts) to Text format so Excel
# does not auto-convert numeric/date-looking strings to numbers/dates.
$ws.Range("A342:A362").NumberFormat = "@"
$ws.Range("C342:C362").NumberFormat = "@"

$ws.Range("A342").Value = "06.04.25"
$ws.Range("B342").Value = "Brindis"
$ws.Range("C342").Value = "103,30"
$ws.Range("D342").Value = "pdf"
$ws.Range("E342").Value = "Facturas/2025/2T/ABRIL/06.04.25 Brindis 103,30.pdf"

$ws.Range("A343").Value = "03.06.20"
$ws.Range("B343").Value = "Hotmart"
$ws.Range("C343").Value = "555"
$ws.Range("D343").Value = "pdf"
$ws.Range("E343").Value = "Facturas/2025/2T/RESAGADAS/03.06.20 Hotmart 555.pdf"

$ws.Range("A344").Value = "04.08.20"
$ws.Range("B344").Value = "Hotmart"
$ws.Range("C344").Value = "49"
$ws.Range("D344").Value = "pdf"
$ws.Range("E344").Value = "Facturas/2025/2T/RESAGADAS/04.08.20 Hotmart 49.pdf"

$ws.Range("A345").Value = "05.01.23"
$ws.Range("B345").Value = "Hotmart"
$ws.Range("C345").Value = "90,21"
$ws.Range("D345").Value = "pdf"
$ws.Range("E345").Value = "Facturas/2025/2T/RESAGADAS/05.01.23 Hotmart 90,21.pdf"

$ws.Range("A346").Value = "05.01.23"
$ws.Range("B346").Value = "Hotmart energia"
$ws.Range("C346").Value = "90,21"
$ws.Range("D346").Value = "pdf"
$ws.Range("E346").Value = "Facturas/2025/2T/RESAGADAS/05.01.23 Hotmart energia 90,21.pdf"

$ws.Range("A347").Value = "05.02.21"
$ws.Range("B347").Value = "Hotmart"
$ws.Range("C347").Value = "97"
$ws.Range("D347").Value = "pdf"
$ws.Range("E347").Value = "Facturas/2025/2T/RESAGADAS/05.02.21 Hotmart 97.pdf"

$ws.Range("A348").Value = "05.12.22"
$ws.Range("B348").Value = "Hotmart"
$ws.Range("C348").Value = "35,09"
$ws.Range("D348").Value = "pdf"
$ws.Range("E348").Value = "Facturas/2025/2T/RESAGADAS/05.12.22 Hotmart 35,09.pdf"

$ws.Range("A349").Value = "06.04.21"
$ws.Range("B349").Value = "Hotmart"
$ws.Range("C349").Value = "457,99"
$ws.Range("D349").Value = "pdf"
$ws.Range("E349").Value = "Facturas/2025/2T/RESAGADAS/06.04.21 Hotmart 457,99.pdf"

$ws.Range("A350").Value = "07.10.24"
$ws.Range("B350").Value = "Atlas FZE"
$ws.Range("C350").Value = "5000"
$ws.Range("D350").Value = "pdf"
$ws.Range("E350").Value = "Facturas/2025/2T/RESAGADAS/07.10.24 Atlas FZE 5000.pdf"

$ws.Range("A351").Value = "07.12.22"
$ws.Range("B351").Value = "Hotmart"
$ws.Range("C351").Value = "26,62"
$ws.Range("D351").Value = "pdf"
$ws.Range("E351").Value = "Facturas/2025/2T/RESAGADAS/07.12.22 Hotmart 26,62.pdf"

$ws.Range("A352").Value = "07.12.22"
$ws.Range("B352").Value = "Hotmart"
$ws.Range("C352").Value = "465,85"
$ws.Range("D352").Value = "pdf"
$ws.Range("E352").Value = "Facturas/2025/2T/RESAGADAS/07.12.22 Hotmart 465,85.pdf"

$ws.Range("A353").Value = "17.02.21"
$ws.Range("B353").Value = "Hotmart"
$ws.Range("C353").Value = "1,11"
$ws.Range("D353").Value = "pdf"
$ws.Range("E353").Value = "Facturas/2025/2T/RESAGADAS/17.02.21 Hotmart 1,11.pdf"

$ws.Range("A354").Value = "19.01.24"
$ws.Range("B354").Value = "Hotmart"
$ws.Range("C354").Value = "601,37"
$ws.Range("D354").Value = "pdf"
$ws.Range("E354").Value = "Facturas/2025/2T/RESAGADAS/19.01.24 Hotmart 601,37.pdf"

$ws.Range("A355").Value = "19.08.22"
$ws.Range("B355").Value = "Hotmart"
$ws.Range("C355").Value = "15,57"
$ws.Range("D355").Value = "pdf"
$ws.Range("E355").Value = "Facturas/2025/2T/RESAGADAS/19.08.22 Hotmart 15,57.pdf"

$ws.Range("A356").Value = "19.10.20"
$ws.Range("B356").Value = "Hotmart"
$ws.Range("C356").Value = "422,29"
$ws.Range("D356").Value = "pdf"
$ws.Range("E356").Value = "Facturas/2025/2T/RESAGADAS/19.10.20 Hotmart 422,29.pdf"

$ws.Range("A357").Value = "19.10.20"
$ws.Range("B357").Value = "Hotmart"
$ws.Range("C357").Value = "97"
$ws.Range("D357").Value = "pdf"
$ws.Range("E357").Value = "Facturas/2025/2T/RESAGADAS/19.10.20 Hotmart 97.pdf"

$ws.Range("A358").Value = "20.01.23"
$ws.Range("B358").Value = "Hotmart"
$ws.Range("C358").Value = "111.32"
$ws.Range("D358").Value = "pdf"
$ws.Range("E358").Value = "Facturas/2025/2T/RESAGADAS/20.01.23 Hotmart 111.32.pdf"

$ws.Range("A359").Value = "23.06.20"
$ws.Range("B359").Value = "Hotmart"
$ws.Range("C359").Value = "3626,37"
$ws.Range("D359").Value = "pdf"
$ws.Range("E359").Value = "Facturas/2025/2T/RESAGADAS/23.06.20 Hotmart 3626,37.pdf"

$ws.Range("A360").Value = "23.12.22"
$ws.Range("B360").Value = "Hotmart"
$ws.Range("C360").Value = "1811.37"
$ws.Range("D360").Value = "pdf"
$ws.Range("E360").Value = "Facturas/2025/2T/RESAGADAS/23.12.22 Hotmart 1811.37.pdf"

$ws.Range("A361").Value = "23.12.22"
$ws.Range("B361").Value = "Hotmart"
$ws.Range("C361").Value = "81,07"
$ws.Range("D361").Value = "pdf"
$ws.Range("E361").Value = "Facturas/2025/2T/RESAGADAS/23.12.22 Hotmart 81,07.pdf"

$ws.Range("A362").Value = "29.10.23"
$ws.Range("B362").Value = "Hotmart"
$ws.Range("C362").Value = "671,55"
$ws.Range("D362").Value = "pdf"
$ws.Range("E362").Value = "Facturas/2025/2T/RESAGADAS/29.10.23 Hotmart 671,55.pdf"
